$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 0. Preserve the original "Additional Site" row (current row 2) formatting so
#    it can be reapplied later, once shifted down to row 6.
# ---------------------------------------------------------------------------
$ws.Range("A2:Q2").Copy() | Out-Null
$ws.Range("A6:Q6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("P2").Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4122) | Out-Null      # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 1. Row 6 -- keep old "Bansal" template formatting, fill with new
#    PITISUTTIHUM Punnee investigator (site-specific Additional-Site view).
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "Principal"
$ws.Range("B6").Value = "0000/0008"
$ws.Range("C6").Value = "'0102"
$ws.Range("D6").Value = "PITISUTTIHUM Punnee"
$ws.Range("G6").Value = "Punnee"
$ws.Range("I6").Value = "PITISUTTIHUM"
$ws.Range("J6").Value = "Apollo"
$ws.Range("K6").Value = "abc street"
$ws.Range("M6").Value = "NY"
$ws.Range("N6").Value = "US"
$ws.Range("O6").Value = 889
$ws.Range("P6").Value = "Spain"
$ws.Range("Q6").Value = "#6789"

# ---------------------------------------------------------------------------
# 2. Row 7 -- blank helper row underneath, same font as "Member Country"
#    label column, no value.
# ---------------------------------------------------------------------------
$ws.Range("D7").Value = ""

# ---------------------------------------------------------------------------
# 3. Row 2 -- first investigator, MACIAS-PARRA Mercedes, with the refreshed
#    "Add Additional Site" view styling (center aligned, no wrap on the
#    name columns, default workbook font throughout).
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Principal"
$ws.Range("B2").Value = "0000/0008"
$ws.Range("C2").Value = "'0102"
$ws.Range("D2").Value = "MACIAS-PARRA, MERCEDES"
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = "MERCEDES"
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = "MACIAS-PARRA"
$ws.Range("J2").Value = "Apollo"
$ws.Range("K2").Value = "abc street"
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = "NY"
$ws.Range("N2").Value = "US"
$ws.Range("O2").Value = 889
$ws.Range("P2").Value = "Spain"
$ws.Range("Q2").Value = "#6789"

# Simplify the formatting of the name columns (D:L) to match the refreshed
# "Add Additional Site" view: centered horizontally, no wrap, no vertical
# centering, default font.
$simpleCols = @("D","E","F","G","H","I","L")
foreach ($col in $simpleCols) {
    $rng = $ws.Range($col + "2")
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 11
    $rng.Font.ColorIndex = -4105
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4160
    $rng.WrapText = $false
}

# ---------------------------------------------------------------------------
# 4. Row 3 -- second investigator, Punnee Pitisuttithum, same styling as row 2.
# ---------------------------------------------------------------------------
$ws.Range("A2:Q2").Copy() | Out-Null
$ws.Range("A3:Q3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A3").Value = "Principal"
$ws.Range("B3").Value = "0000/0008"
$ws.Range("C3").Value = "'0102"
$ws.Range("D3").Value = "Punnee Pitisuttithum"
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = "Punnee"
$ws.Range("H3").Value = ""
$ws.Range("I3").Value = "Pitisuttithum"
$ws.Range("J3").Value = "Apollo"
$ws.Range("K3").Value = "abc street"
$ws.Range("L3").Value = ""
$ws.Range("M3").Value = "NY"
$ws.Range("N3").Value = "US"
$ws.Range("O3").Value = 889
$ws.Range("P3").Value = "Spain"
$ws.Range("Q3").Value = "#6789"

# ---------------------------------------------------------------------------
# 5. Rows 4 & 5 -- blank filler rows carrying the same (new) formatting.
# ---------------------------------------------------------------------------
$ws.Range("A3:Q3").Copy() | Out-Null
$ws.Range("A4:Q4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A5:Q5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 6. Sheet view / dimension bookkeeping to match the refreshed view.
# ---------------------------------------------------------------------------
$ws.Range("I3").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# Column width refresh for Middle Name / Last Name (content shrank).
$ws.Columns(8).EntireColumn.AutoFit() | Out-Null
$ws.Columns(9).EntireColumn.AutoFit() | Out-Null
